$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from row 168 down into rows 169 and 170
$ws.Range("A168:I168").Copy()
$ws.Range("A169:I170").PasteSpecial(-4122)
$ws.Range("A169:I170").PasteSpecial(-4142)
$ws.Rows.Item(169).RowHeight = 13.4
$ws.Rows.Item(170).RowHeight = 13.4

# Row 169 - LeetCode #175 "Combine Two Tables" (SQL)
$ws.Cells.Item(169, 1).Value = 175
$ws.Cells.Item(169, 2).Value = "简单"
$ws.Cells.Item(169, 3).Value = "组合两个表"
$ws.Cells.Item(169, 4).Value = "combine-two-tables"
$ws.Cells.Item(169, 5).Value = "200221-1.sql"
$ws.Cells.Item(169, 6).Value = 176
$ws.Cells.Item(169, 7).Value = 0
$ws.Cells.Item(169, 8).Value = "41.84%"
$ws.Cells.Item(169, 9).Value = 48466380

# Row 170 - LeetCode #176 "Second Highest Salary" (SQL)
$ws.Cells.Item(170, 1).Value = 176
$ws.Cells.Item(170, 2).Value = "简单"
$ws.Cells.Item(170, 3).Value = "第二高的薪水"
$ws.Cells.Item(170, 4).Value = "second-highest-salary"
$ws.Cells.Item(170, 5).Value = "200221-1.cpp"
$ws.Cells.Item(170, 6).Value = 181
$ws.Cells.Item(170, 7).Value = 0
$ws.Cells.Item(170, 8).Value = "10.64%"
$ws.Cells.Item(170, 9).Value = 48466681

$ws.Range("I171").Select()


